# Refresh scraped market-price columns (H:N) across the Leve tables.
# Generated from the authoritative before/after cell diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1306.5
$ws.Range("I2").Value = 1399.3077
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1399.3077
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -1286.3077
$ws.Range("N2").Value = -326
$ws.Range("H28").Value = 573.4211
$ws.Range("I28").Value = 630.1875
$ws.Range("J28").Value = 270.66666
$ws.Range("K28").Value = 630.1875
$ws.Range("L28").Value = 270.66666
$ws.Range("M28").Value = -145.1875
$ws.Range("N28").Value = -1240.66666
$ws.Range("H40").Value = 5380.875
$ws.Range("J40").Value = 7442.625
$ws.Range("L40").Value = 7442.625
$ws.Range("N40").Value = -7792.625
$ws.Range("H132").Value = 3567.52
$ws.Range("I132").Value = 916.3415
$ws.Range("J132").Value = 15645.111
$ws.Range("K132").Value = 2749.0245
$ws.Range("L132").Value = 46935.333
$ws.Range("M132").Value = -219.0245
$ws.Range("N132").Value = -51995.333
$ws.Range("H135").Value = 227.5
$ws.Range("I135").Value = 227.5
$ws.Range("K135").Value = 2047.5
$ws.Range("M135").Value = 487.5
$ws.Range("H137").Value = 1981691.4
$ws.Range("I137").Value = 2025422.8
$ws.Range("K137").Value = 6076268.4
$ws.Range("M137").Value = -6073718.4
$ws.Range("H138").Value = 1659.3695
$ws.Range("I138").Value = 677.62067
$ws.Range("J138").Value = 3334.1177
$ws.Range("K138").Value = 2032.86201
$ws.Range("L138").Value = 10002.3531
$ws.Range("M138").Value = 3107.13799
$ws.Range("N138").Value = -20282.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1638652.4
$ws.Range("I32").Value = 1705332
$ws.Range("K32").Value = 1705332
$ws.Range("M32").Value = -1705045
$ws.Range("H110").Value = 2232.5833
$ws.Range("I110").Value = 2798.7144
$ws.Range("J110").Value = 1440
$ws.Range("K110").Value = 2798.7144
$ws.Range("L110").Value = 1440
$ws.Range("M110").Value = -753.7143999999998
$ws.Range("N110").Value = -5530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1980370.6
$ws.Range("I31").Value = 1395.0769
$ws.Range("J31").Value = 2759967
$ws.Range("K31").Value = 1395.0769
$ws.Range("L31").Value = 2759967
$ws.Range("M31").Value = -1100.0769
$ws.Range("N31").Value = -2760557
$ws.Range("H34").Value = 1980370.6
$ws.Range("I34").Value = 1395.0769
$ws.Range("J34").Value = 2759967
$ws.Range("K34").Value = 1395.0769
$ws.Range("L34").Value = 2759967
$ws.Range("M34").Value = -1193.0769
$ws.Range("N34").Value = -2760371
$ws.Range("H62").Value = 2781283.8
$ws.Range("I62").Value = 5557989.5
$ws.Range("K62").Value = 5557989.5
$ws.Range("M62").Value = -5557365.5
$ws.Range("H65").Value = 2781283.8
$ws.Range("I65").Value = 5557989.5
$ws.Range("K65").Value = 27789947.5
$ws.Range("M65").Value = -27786827.5
$ws.Range("H107").Value = 834.7778
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 1002.6
$ws.Range("K107").Value = 625
$ws.Range("L107").Value = 1002.6
$ws.Range("M107").Value = 1295
$ws.Range("N107").Value = -4842.6
$ws.Range("H134").Value = 31252016
$ws.Range("I134").Value = 71429656
$ws.Range("J134").Value = 2739.7778
$ws.Range("K134").Value = 214288968
$ws.Range("L134").Value = 8219.3334
$ws.Range("M134").Value = -214286433
$ws.Range("N134").Value = -13289.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.958336
$ws.Range("I2").Value = 99.57143000000001
$ws.Range("J2").Value = 69.3
$ws.Range("K2").Value = 597.42858
$ws.Range("L2").Value = 415.8
$ws.Range("M2").Value = -484.42858
$ws.Range("N2").Value = -641.8
$ws.Range("H7").Value = 406.53845
$ws.Range("I7").Value = 99.8
$ws.Range("J7").Value = 598.25
$ws.Range("K7").Value = 299.4
$ws.Range("L7").Value = 1794.75
$ws.Range("M7").Value = -187.4
$ws.Range("N7").Value = -2018.75
$ws.Range("H17").Value = 10001
$ws.Range("I17").Value = 10001
$ws.Range("K17").Value = 30003
$ws.Range("M17").Value = -29834
$ws.Range("H23").Value = 595.9
$ws.Range("I23").Value = 61
$ws.Range("J23").Value = 655.3333
$ws.Range("K23").Value = 183
$ws.Range("L23").Value = 1965.9999
$ws.Range("M23").Value = 52
$ws.Range("N23").Value = -2435.9999
$ws.Range("H33").Value = 329.73914
$ws.Range("I33").Value = 250.15384
$ws.Range("J33").Value = 433.2
$ws.Range("K33").Value = 1500.92304
$ws.Range("L33").Value = 2599.2
$ws.Range("M33").Value = -1217.92304
$ws.Range("N33").Value = -3165.2
$ws.Range("H34").Value = 448.47827
$ws.Range("I34").Value = 67.083336
$ws.Range("J34").Value = 864.5454999999999
$ws.Range("K34").Value = 201.250008
$ws.Range("L34").Value = 2593.6365
$ws.Range("M34").Value = -117.250008
$ws.Range("N34").Value = -2761.6365
$ws.Range("H39").Value = 4427.273
$ws.Range("J39").Value = 4427.273
$ws.Range("L39").Value = 13281.819
$ws.Range("N39").Value = -13869.819
$ws.Range("H55").Value = 2811.4707
$ws.Range("J55").Value = 2962.1875
$ws.Range("L55").Value = 8886.5625
$ws.Range("N55").Value = -9240.5625
$ws.Range("H86").Value = 707.0714
$ws.Range("I86").Value = 700
$ws.Range("J86").Value = 733
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 2199
$ws.Range("M86").Value = -914
$ws.Range("N86").Value = -4571
$ws.Range("H89").Value = 707.0714
$ws.Range("I89").Value = 700
$ws.Range("J89").Value = 733
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 6597
$ws.Range("M89").Value = -372
$ws.Range("N89").Value = -18453
$ws.Range("H118").Value = 2476.6428
$ws.Range("I118").Value = 510.7143
$ws.Range("J118").Value = 4442.5713
$ws.Range("K118").Value = 1532.1429
$ws.Range("L118").Value = 13327.7139
$ws.Range("M118").Value = -289.1428999999998
$ws.Range("N118").Value = -15813.7139
$ws.Range("H122").Value = 717.4286
$ws.Range("I122").Value = 309.7143
$ws.Range("J122").Value = 1125.1428
$ws.Range("K122").Value = 2787.4287
$ws.Range("L122").Value = 10126.2852
$ws.Range("M122").Value = -337.4286999999999
$ws.Range("N122").Value = -15026.2852
$ws.Range("H125").Value = 4635.8335
$ws.Range("I125").Value = 630
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 1890
$ws.Range("L125").Value = 15000
$ws.Range("M125").Value = 3030
$ws.Range("N125").Value = -24840
$ws.Range("H137").Value = 8944.333000000001
$ws.Range("I137").Value = 1616.6666
$ws.Range("J137").Value = 13829.444
$ws.Range("K137").Value = 4849.9998
$ws.Range("L137").Value = 41488.33199999999
$ws.Range("M137").Value = 250.0002000000004
$ws.Range("N137").Value = -51688.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1260
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 2510
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 2510
$ws.Range("M2").Value = 103
$ws.Range("N2").Value = -2736
$ws.Range("H98").Value = 14529.615
$ws.Range("J98").Value = 14529.615
$ws.Range("L98").Value = 14529.615
$ws.Range("N98").Value = -20519.615
$ws.Range("H101").Value = 46666.668
$ws.Range("J101").Value = 46666.668
$ws.Range("L101").Value = 46666.668
$ws.Range("N101").Value = -53156.668
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 1925862.8
$ws.Range("I132").Value = 2978591.5
$ws.Range("J132").Value = 3488.5217
$ws.Range("K132").Value = 8935774.5
$ws.Range("L132").Value = 10465.5651
$ws.Range("M132").Value = -8933244.5
$ws.Range("N132").Value = -15525.5651

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 4983
$ws.Range("I33").Value = 4978.75
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 4978.75
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -4688.75
$ws.Range("N33").Value = -5580
$ws.Range("H132").Value = 17367566
$ws.Range("I132").Value = 14288075
$ws.Range("J132").Value = 25658506
$ws.Range("K132").Value = 42864225
$ws.Range("L132").Value = 76975518
$ws.Range("M132").Value = -42861695
$ws.Range("N132").Value = -76980578
$ws.Range("H136").Value = 3920.465
$ws.Range("I136").Value = 2171.7727
$ws.Range("J136").Value = 5752.4287
$ws.Range("K136").Value = 6515.3181
$ws.Range("L136").Value = 17257.2861
$ws.Range("M136").Value = -3965.3181
$ws.Range("N136").Value = -22357.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1698.6666
$ws.Range("I126").Value = 1053.8667
$ws.Range("J126").Value = 3310.6667
$ws.Range("K126").Value = 3161.6001
$ws.Range("L126").Value = 9932.000100000001
$ws.Range("M126").Value = -691.6001000000001
$ws.Range("N126").Value = -14872.0001
$ws.Range("H132").Value = 15874720
$ws.Range("I132").Value = 20409328
$ws.Range("J132").Value = 3592
$ws.Range("K132").Value = 61227984
$ws.Range("L132").Value = 10776
$ws.Range("M132").Value = -61225454
$ws.Range("N132").Value = -15836
